# Added new login backup codes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift the newest backup codes up into A2:A4 (a new code was inserted at the
# top of the list, pushing the old ones down and bumping the last code off
# the visible/used block at A7).
$ws.Range("A2").Value = "51K0DF5KCN34"
$ws.Range("A3").Value = "HQT8HMXSF63S"
$ws.Range("A4").Value = "5A41AVCX9PFR"

# The code that used to live at A7 has now been consumed/removed.
$ws.Range("A7").ClearContents()

# Update the active selection to reflect where the user left off.
$ws.Range("A4").Select()
